# Update course Excel file: replace the "SHELDON SCHOOL OF HOSPITALITY"
# department label with "Packages" (rows 2-5) / "Hospitality" (row 6),
# and adjust the related row heights / selection to match the resave.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# Column C ("department") text changes
$ws.Range("C2").Value = "Packages"
$ws.Range("C3").Value = "Packages"
$ws.Range("C4").Value = "Packages"
$ws.Range("C5").Value = "Packages"
$ws.Range("C6").Value = "Hospitality"

# Row heights now uniform at 42.75 (rows 2 & 3 were 45, rows 4-6 were 60)
$ws.Rows.Item(2).RowHeight = 42.75
$ws.Rows.Item(3).RowHeight = 42.75
$ws.Rows.Item(4).RowHeight = 42.75
$ws.Rows.Item(5).RowHeight = 42.75
$ws.Rows.Item(6).RowHeight = 42.75

# Active selection moved to C6 (last cell edited) before save
[void]$ws.Range("C6").Select()
